$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old F1 value (column no longer used)
$ws.Range("F1").ClearContents()

# New header/filter row values
$ws.Range("A1").Value = "Купить"
$ws.Range("B1").Value = "Комната"
$ws.Range("C1").Value = 1000
$ws.Range("D1").Value = 1000000000
$ws.Range("E1").Value = "Воронеж"

# Column A gets a custom width, matching the new layout
$ws.Columns.Item(1).ColumnWidth = 14.83

# Move the active selection to the new last used cell
$ws.Range("E1").Select() | Out-Null
